$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest cryptos snapshot.
# D-column values are plain text (e.g. "27.138.74"), and some of the new readings
# (e.g. "216.18") look like ordinary decimals, so Excel would otherwise auto-convert
# them to Number on assignment. Force Text number-format for the write, then clear the
# format back off so the cell keeps its original (default) style -- only its value
# changes, matching the source refresh.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue "D2" '27.138.74'
$ws.Range("E2").Value = '  -0.30%  '
Set-TextValue "D3" '1.631.65'
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("E4").Value = '  -0.03%  '
Set-TextValue "D5" '216.18'
$ws.Range("E5").Value = '  -0.96%  '
Set-TextValue "D6" '0.517'
$ws.Range("E6").Value = '  +1.22%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -1.26%  '
$ws.Range("E9").Value = '  -0.66%  '
$ws.Range("E10").Value = '  -0.59%  '
$ws.Range("E11").Value = '  +0.24%  '
Set-TextValue "D12" '1.859.84'
$ws.Range("E12").Value = '  -0.97%  '
Set-TextValue "D13" '1.655.51'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("E15").Value = '  +0.46%  '
Set-TextValue "D16" '65.67'
$ws.Range("E16").Value = '  -3.03%  '
Set-TextValue "D17" '27.110.62'
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("E18").Value = '  -1.07%  '
Set-TextValue "D19" '214.48'
$ws.Range("E19").Value = '  -2.77%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("E21").Value = '  +0.91%  '
$ws.Range("E22").Value = '  -1.12%  '
$ws.Range("E23").Value = '  -0.61%  '
$ws.Range("E24").Value = '  -1.17%  '
Set-TextValue "D25" '147.21'
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  -0.50%  '
$ws.Range("E28").Value = '  -1.24%  '
$ws.Range("E29").Value = '  -1.21%  '
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("E33").Value = '  -1.01%  '
Set-TextValue "D34" '1.309.88'
$ws.Range("E34").Value = '  +2.82%  '
$ws.Range("E35").Value = '  -1.38%  '
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("E37").Value = '  -1.51%  '
Set-TextValue "D38" '0.543'
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("E42").Value = '  -0.68%  '
$ws.Range("E43").Value = '  -1.96%  '
Set-TextValue "D44" '1.769.19'
Set-TextValue "D45" '62.18'
$ws.Range("E45").Value = '  -1.08%  '
Set-TextValue "D46" '90.58'
$ws.Range("E46").Value = '  -1.92%  '
$ws.Range("E47").Value = '  +0.32%  '
$ws.Range("E48").Value = '  +0.18%  '
Set-TextValue "D49" '0.818'
$ws.Range("E49").Value = '  +21.25%  '
$ws.Range("E50").Value = '  -0.07%  '
Set-TextValue "D51" '7.58'
$ws.Range("E51").Value = '  -1.17%  '
